$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.401.35"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").Value = "2.600.82"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.34"
$ws.Range("E5").Value = "  +3.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.60"
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").Value = "2.611.46"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("E11").Value = "  +3.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("E12").Value = "  +2.48%  "
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").Value = "3.060.67"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").Value = "59.334.74"
$ws.Range("E15").Value = "  +2.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.54"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "2.629.48"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "345.10"
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.33"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.10"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.35"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.02"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  +1.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.20"
$ws.Range("E28").Value = "  +3.52%  "
$ws.Range("D29").Value = "0.0₃0742"
$ws.Range("E29").Value = "  +4.86%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +5.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.83"
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.84"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.13"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.99"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.842"
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.47"
$ws.Range("E39").Value = "  +3.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.836"
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "277.07"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.599"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0960"
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0523"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("D48").Value = "1.949.97"
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0223"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.34"
$ws.Range("E50").Value = "  +4.20%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.50"
$ws.Range("E51").Value = "  -0.14%  "
